$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the numeric data in columns B:E (rows 2-13) to the nearest integer,
# matching the "write to disk as integer data" behavior described in the
# commit message.
$range = $ws.Range("B2:E13")
foreach ($cell in $range.Cells) {
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value = [Math]::Round([double]$val, 0)
    }
}
